$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster label changes from "Neutrophils" to "ECs" because of shared-string
# reordering in the source diff; recompute the dependent numeric columns.
$ws.Range("A2").Value = "ECs"
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.0565965
$ws.Range("H2").Value = 0.113193
$ws.Range("O2").Value = 0.2281429585585933
$ws.Range("P2").Value = 0.2281429585585933
$ws.Range("Q2").Value = 0.0026797499475
$ws.Range("R2").Value = 0.016078499685
$ws.Range("S2").Value = 0.2281429585585933
$ws.Range("T2").Value = 0.2281429585585933

# New row 3: ECs -> Pnoc/Oprl1 -> Neutrophils
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pnoc"
$ws.Range("C3").Value = "Oprl1"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.0565965
$ws.Range("H3").Value = 0.113193
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1601896666666667
$ws.Range("N3").Value = 0.480569
$ws.Range("O3").Value = 0.7718570414414068
$ws.Range("P3").Value = 0.7718570414414068
$ws.Range("Q3").Value = 0.009066174469500001
$ws.Range("R3").Value = 0.05439704681700001
$ws.Range("S3").Value = 0.7718570414414068
$ws.Range("T3").Value = 0.7718570414414068
